# RBA v2.5 - Atualizacao da Tela
#
# Replaces the placeholder "TRE/TERE/Tre/tre" style tokens with
# "QWER/Qwer/Qewr/qwer" tokens both in the body of the letter and in the
# page header, one occurrence at a time (several runs in the header
# share identical text, e.g. multiple runs containing "Tre", but each
# must get a different replacement value, so we replace them in
# document order using repeated Find.Execute calls on the same Range,
# which always finds the next occurrence after the previous match).
#
# Each replacement is wrapped in its own TrackRevisions / AcceptAllRevisions
# cycle. Word's plain Find&Replace silently coalesces the edited run into
# an adjacent run that happens to carry identical direct formatting
# (common in this document, since neighboring runs only differ by
# w:rsid* bookkeeping attributes), which would delete run/paragraph
# structure that the target revision keeps intact. Performing the edit
# as a tracked insertion+deletion and accepting it immediately keeps
# the run boundaries untouched while still leaving plain (non-tracked)
# text behind.

$d = $word.ActiveDocument

function Replace-OneOccurrence($rng, $find, $repl) {
    $d.TrackRevisions = $true
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 1)
    $d.TrackRevisions = $false
    $d.AcceptAllRevisions()
}

# --- Main document body ---------------------------------------------------
# "A TERE," -> "A QWER," (bold run inside the body paragraph)
$body = $d.Content
Replace-OneOccurrence $body "TERE" "QWER"

# --- Page header ------------------------------------------------------------
$hdr = $d.Sections.First.Headers.Item(1)
$hrng = $hdr.Range

# "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
Replace-OneOccurrence $hrng "TRE" "QWER"
# "TERE - DEP." -> "QWER - DEP."
Replace-OneOccurrence $hrng "TERE" "QWER"
# "Tre, nº Tre - Tre - Tre - Tre" -> "Qwer, nº Qwer - Qewr - Qewr - Qwer"
Replace-OneOccurrence $hrng "Tre" "Qwer"
Replace-OneOccurrence $hrng "Tre" "Qwer"
Replace-OneOccurrence $hrng "Tre" "Qewr"
Replace-OneOccurrence $hrng "Tre" "Qewr"
Replace-OneOccurrence $hrng "Tre" "Qwer"
# "CEP: tre" -> "CEP: qwer"
Replace-OneOccurrence $hrng "tre" "qwer"
# "Tel: tre" -> "Tel: qwer"
Replace-OneOccurrence $hrng "tre" "qwer"
# "Email: tre" -> "Email: qwer"
Replace-OneOccurrence $hrng "tre" "qwer"
